# The "Chrome/Safari browser" post occupied row 319 and was removed from
# the source data. Deleting the entire row shifts every row below it up
# by one (old row 320 -> new row 319, ..., old row 454 -> new row 453),
# which also shrinks the sheet's used range from A1:C454 to A1:C453 -
# exactly matching the target diff (no cell content below row 319 actually
# changes, only its row number).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("319").Delete()
